$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF"), matching the style of the existing header row (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate I0/IF numeric values for data rows 2-49
$data = @(
  @(2, 9, 9),
  @(3, 7, 7),
  @(4, 8, 8),
  @(5, 8, 8),
  @(6, 7, 7),
  @(7, 8, 9),
  @(8, 6, 7),
  @(9, 6, 7),
  @(10, 5, 5),
  @(11, 8, 8),
  @(12, 7, 7),
  @(13, 7, 7),
  @(14, 6, 7),
  @(15, 6, 6),
  @(16, 10, 10),
  @(17, 6, 6),
  @(18, 9, 9),
  @(19, 5, 5),
  @(20, 5, 6),
  @(21, 6, 6),
  @(22, 7, 8),
  @(23, 7, 7),
  @(24, 7, 7),
  @(25, 6, 7),
  @(26, 6, 7),
  @(27, 6, 6),
  @(28, 6, 6),
  @(29, 5, 5),
  @(30, 7, 8),
  @(31, 7, 7),
  @(32, 6, 6),
  @(33, 8, 8),
  @(34, 12, 12),
  @(35, 8, 8),
  @(36, 7, 8),
  @(37, 9, 9),
  @(38, 7, 8),
  @(39, 10, 10),
  @(40, 5, 6),
  @(41, 8, 9),
  @(42, 9, 9),
  @(43, 8, 8),
  @(44, 7, 7),
  @(45, 6, 6),
  @(46, 6, 6),
  @(47, 7, 7),
  @(48, 3, 3),
  @(49, 3, 3)
)

foreach ($row in $data) {
  $r = $row[0]
  $i0 = $row[1]
  $ifv = $row[2]
  $ws.Cells.Item($r, 9).Value = $i0
  $ws.Cells.Item($r, 10).Value = $ifv
}

Write-Output "done"
